$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "GRT-USD"
$ws.Range("A19").Value = "BSCX-USD"
